$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the UML diagram text so the width/height setters reflect that they
# now return a bool (indicating success) instead of void.
$ws.Range("B19").Value = " +setWidth(width: double): bool"
$ws.Range("B21").Value = " +setHeight(height: double): bool"

# Update the saved selection to match the author's final cursor position.
$ws.Range("D21").Select()
